# Product Backlog.xlsx edit
# Mark the first four backlog items ("to move through the world",
# "see a description of current location", "see action prompts",
# "see player details") as done by checking the box in column A
# (rows 2-5), mirroring the author's in-cell checkbox toggle, and
# leave the selection on A5 where the cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").Value = $true

$ws.Range("A5").Select() | Out-Null
